$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Image_Path for "Pain d'épices" (row 5, column D) from .jpg to .png
$ws.Range("D5").Value = "https://raw.githubusercontent.com/AlDenervaud/champdupuits/refs/heads/main/data/images/apiculture/pain_epices.png"

# Update the active selection to D6 as recorded in the saved workbook view
$ws.Range("D6").Select()
